$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.656.38"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.801.15"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").Value = "'1.005"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "'327.77"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("D6").Value = "'1.002"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4368"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  -4.78%  "
$ws.Range("D8").Value = "'0.3745"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("D9").Value = "'46.10"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "'0.07599"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "'1.134"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").Value = "'22.62"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'1.004"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'6.230"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "'7.484"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "1.803.35"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.00001086"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "'0.06690"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'80.87"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "'1.002"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "'17.54"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "'6.229"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "28.662.11"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "'11.71"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Value = "'2.436"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "'20.43"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").Value = "'154.23"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").Value = "'2.335"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").Value = "2.010.94"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'1.299"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").Value = "'130.64"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").Value = "'3.970"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").Value = "'5.788"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("D34").Value = "'0.09173"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").Value = "'0.2228"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "'12.11"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Value = "'0.06276"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'0.02311"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("D39").Value = "'5.167"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "'0.6569"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("D41").Value = "'1.196"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").Value = "'8.031"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").Value = "'1.427"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  -4.43%  "
$ws.Range("D44").Value = "'1.001"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'13.87"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "'0.6066"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "'127.38"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").Value = "'2.014"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "'0.07023"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  -1.17%  "

# Row 51: EOS -> Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'77.66"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
